$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# Fix D15/D16 to be numeric values instead of text
$ws.Range("D15").Value = 524494
$ws.Range("D16").Value = 509480

# Add new row 17 - CUMMINSIND
$ws.Range("A17").Value = 1
$ws.Range("B17").Value = "CUMMINSIND"
$ws.Range("C17").Value = "Cummins India Limited"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "500480"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = 1.02
$ws.Range("F17").Value = 3712.5
$ws.Range("G17").Value = 978774
$ws.Range("H17").Value = "day"
$ws.Range("I17").Value = "13/06/2024 10:33:04"

# Add new row 18 - UBL
$ws.Range("A18").Value = 2
$ws.Range("B18").Value = "UBL"
$ws.Range("C18").Value = "United Breweries Limited"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "532478"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = -0.45
$ws.Range("F18").Value = 2120
$ws.Range("G18").Value = 84876
$ws.Range("H18").Value = "day"
$ws.Range("I18").Value = "13/06/2024 10:33:04"
